$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.605.29'
$ws.Range('E2').Value = '  +1.18%  '
$ws.Range('D3').Value = '1.798.76'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('E4').Value = '  -0.09%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '227.28'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('E6').Value = '  +2.02%  '
$ws.Range('E7').Value = '  -0.09%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '32.86'
$c.ClearFormats()
$ws.Range('E8').Value = '  +3.17%  '
$ws.Range('E9').Value = '  +1.91%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0698'
$c.ClearFormats()
$ws.Range('E10').Value = '  +1.04%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0952'
$c.ClearFormats()
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').Value = '2.058.14'
$ws.Range('E12').Value = '  +0.94%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '11.12'
$c.ClearFormats()
$ws.Range('E13').Value = '  +0.95%  '
$ws.Range('D14').Value = '1.805.31'
$ws.Range('E14').Value = '  +1.41%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.640'
$c.ClearFormats()
$ws.Range('E15').Value = '  +2.62%  '
$ws.Range('D16').Value = '34.582.79'
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('E17').Value = '  +3.05%  '
$ws.Range('E18').Value = '  +1.53%  '
$ws.Range('D19').Value = '0.0₃0805'
$ws.Range('E19').Value = '  +1.23%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '247.58'
$c.ClearFormats()
$ws.Range('E20').Value = '  +0.17%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '11.30'
$c.ClearFormats()
$ws.Range('E21').Value = '  +2.77%  '
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('E23').Value = '  +1.84%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '167.86'
$c.ClearFormats()
$ws.Range('E24').Value = '  +3.28%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.08'
$c.ClearFormats()
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('E26').Value = '  +1.63%  '
$ws.Range('E27').Value = '  +1.83%  '
$ws.Range('E28').Value = '  +2.32%  '
$ws.Range('E29').Value = '  -0.12%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '4.10'
$c.ClearFormats()
$ws.Range('E30').Value = '  +11.18%  '
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('E32').Value = '  +1.09%  '
$ws.Range('E33').Value = '  +1.95%  '
$ws.Range('E34').Value = '  +2.74%  '
$ws.Range('D35').Value = '1.431.38'
$ws.Range('E35').Value = '  -0.74%  '
$ws.Range('E36').Value = '  +6.81%  '
$ws.Range('E37').Value = '  +2.81%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '1.07'
$c.ClearFormats()
$ws.Range('E38').Value = '  +2.46%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.0192'
$c.ClearFormats()
$ws.Range('E39').Value = '  +0.62%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '85.62'
$c.ClearFormats()
$ws.Range('E40').Value = '  +6.69%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.42'
$c.ClearFormats()
$ws.Range('E41').Value = '  +0.85%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.936'
$c.ClearFormats()
$ws.Range('E42').Value = '  +1.17%  '
$ws.Range('E43').Value = '  +3.33%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '13.74'
$c.ClearFormats()
$ws.Range('E44').Value = '  +0.49%  '
$ws.Range('E45').Value = '  +3.59%  '
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('D48').Value = '1.957.36'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '106.18'
$c.ClearFormats()
$ws.Range('E49').Value = '  +1.48%  '
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('D51').Value = '0.0₆0128'
$ws.Range('E51').Value = '  -6.81%  '
